$p = $ppt.ActivePresentation

# Slide 1: consolidate "First" + " " + "slide" runs into a single run "First slide"
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "ZZZZZ"
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "First slide"

# Slide 3: consolidate "Third" + " " + "slide" runs into a single run "Third slide"
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "ZZZZZ"
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Third slide"
